$wb = $excel.ActiveWorkbook

# --- Revert the Metadata sheet values back to the pre-2.0 release info ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"

# --- Remove the 13th "Include from FSIII 12" sheet that was added by the reverted commit ---
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Include from FSIII 12").Delete()
$excel.DisplayAlerts = $true
